$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.443.22"
$ws.Range("E2").Value = "  +5.27%  "
$ws.Range("D3").Value = "1.821.34"
$ws.Range("E3").Value = "  +6.00%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'345.76"
$ws.Range("E5").Value = "  +4.58%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "'0.3832"
$ws.Range("E7").Value = "  +3.69%  "
$ws.Range("D8").Value = "'0.3531"
$ws.Range("E8").Value = "  +5.92%  "
$ws.Range("D9").Value = "'49.51"
$ws.Range("E10").Value = "  +4.69%  "
$ws.Range("D11").Value = "'0.07806"
$ws.Range("E11").Value = "  +4.35%  "
$ws.Range("D12").Value = "'1.004"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "'22.24"
$ws.Range("E13").Value = "  +11.03%  "
$ws.Range("D14").Value = "'6.664"
$ws.Range("E14").Value = "  +6.28%  "
$ws.Range("D15").Value = "'7.277"
$ws.Range("E15").Value = "  +5.12%  "
$ws.Range("D16").Value = "1.823.38"
$ws.Range("E16").Value = "  +7.00%  "
$ws.Range("D17").Value = "'0.00001130"
$ws.Range("E17").Value = "  +4.84%  "
$ws.Range("D18").Value = "'0.06756"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("D19").Value = "'86.72"
$ws.Range("E19").Value = "  +5.68%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  +8.67%  "
$ws.Range("D22").Value = "'6.582"
$ws.Range("E22").Value = "  +8.39%  "
$ws.Range("D23").Value = "'13.32"
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("D24").Value = "27.469.99"
$ws.Range("E24").Value = "  +5.66%  "
$ws.Range("D25").Value = "'2.461"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").Value = "'2.704"
$ws.Range("E26").Value = "  +8.75%  "
$ws.Range("D27").Value = "'22.31"
$ws.Range("E27").Value = "  +15.63%  "
$ws.Range("D28").Value = "'1.512"
$ws.Range("E28").Value = "  +13.83%  "
$ws.Range("D29").Value = "'154.09"
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("D30").Value = "2.026.43"
$ws.Range("E30").Value = "  +6.96%  "
$ws.Range("D31").Value = "'137.17"
$ws.Range("E31").Value = "  +6.46%  "
$ws.Range("D32").Value = "'6.407"
$ws.Range("E32").Value = "  +7.29%  "
$ws.Range("D33").Value = "'4.076"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").Value = "'14.16"
$ws.Range("E34").Value = "  +9.45%  "
$ws.Range("D35").Value = "'0.08792"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").Value = "'1.705"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "'5.685"
$ws.Range("E37").Value = "  +5.87%  "
$ws.Range("D38").Value = "'0.7096"
$ws.Range("E38").Value = "  +14.92%  "
$ws.Range("D39").Value = "'0.2291"
$ws.Range("E39").Value = "  +7.30%  "
$ws.Range("E40").Value = "  +5.76%  "
$ws.Range("D41").Value = "'0.02428"
$ws.Range("E41").Value = "  +6.08%  "
$ws.Range("D42").Value = "'9.058"
$ws.Range("E42").Value = "  +6.36%  "
$ws.Range("D43").Value = "'1.301"
$ws.Range("E43").Value = "  +1.79%  "
$ws.Range("D44").Value = "'14.95"
$ws.Range("E44").Value = "  +3.06%  "
$ws.Range("D45").Value = "'0.6625"
$ws.Range("E45").Value = "  +12.51%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "'4.046"
$ws.Range("E47").Value = "  +5.45%  "
$ws.Range("D48").Value = "'2.198"
$ws.Range("E48").Value = "  +9.24%  "
$ws.Range("D49").Value = "'133.29"
$ws.Range("E49").Value = "  +4.60%  "
$ws.Range("D50").Value = "'0.07368"
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("D51").Value = "'81.14"
$ws.Range("E51").Value = "  +5.24%  "
